$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.278.69'
$ws.Range("E2").Value = '  +1.16%  '

$ws.Range("D3").Value = '3.454.07'
$ws.Range("E3").Value = '  +2.36%  '

$ws.Range("D4").NumberFormat = '@'
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '409.73'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +1.55%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '130.56'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +1.12%  '

$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.611'
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  +2.97%  '

$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.747'
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +12.13%  '

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.149'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  +25.32%  '

$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '43.03'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  +4.18%  '

$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '0.141'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  +0.07%  '

$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '0.0000207'
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  +64.07%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '8.73'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  +4.67%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '20.17'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  +3.78%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.460.44'
$ws.Range("E16").Value = '  +2.59%  '

$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '11.76'
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  +5.27%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '62.249.33'
$ws.Range("E18").Value = '  +1.15%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '1.03'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +1.72%  '

$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '3.19'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -0.99%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '332.26'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +9.08%  '

$ws.Range("B22").Value = 'Litecoin'
$ws.Range("C22").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '85.49'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  +3.02%  '

$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '12.99'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  +2.20%  '

$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '3.21'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  +2.99%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '30.85'
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  +5.40%  '

$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '4.76'
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  -0.41%  '

$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '8.26'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -0.26%  '

$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '7.83'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  +4.17%  '

$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '46.13'
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  +11.99%  '

$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '2.75'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  +10.42%  '

$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '0.117'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  +1.44%  '

$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '0.171'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  -0.73%  '

$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '11.76'
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  +3.70%  '

$ws.Range("E34").Value = '  -0.04%  '

$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '0.0487'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  +1.58%  '

$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '51.79'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  -0.27%  '

$ws.Range("E37").Value = '  -0.06%  '

$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '3.38'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  +0.33%  '

$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '2.94'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  -0.23%  '

$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '0.317'
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  +7.13%  '

$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '143.57'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  +4.07%  '

$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '0.129'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  +4.10%  '

$ws.Range("E43").Value = '  +1.11%  '

$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '3.98'
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  +1.83%  '

$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '17.02'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  +2.17%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '2.33'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +5.12%  '

$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '21.74'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  +2.20%  '

$ws.Range("D48").Value = '2.120.18'
$ws.Range("E48").Value = '  +0.51%  '

$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '1.97'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  +4.13%  '

$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '2.29'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  +0.33%  '

$ws.Range("B51").Value = 'Fetch.AI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '1.71'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +7.61%  '
